$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Clear the header-number row (row 4, A4:T4) that was removed in the edit.
# The row had no formatting of its own, so once its contents are gone
# Excel drops the now-empty <row> element entirely on save.
[void]$ws.Range("A4:T4").ClearContents()

# Scroll sheet view back to top-left (removes topLeftCell="L1").
[void]$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

# Restore the active selection to A4:T4, matching the saved sheetView.
[void]$ws.Range("A4:T4").Select()
